$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "woshikuzzaman.cto@gmail.com"
$ws.Range("B3").Value = "703598238217-r6h9ncdqfk52ai610avro7lrk3ecoen1.apps.googleusercontent.com"
$ws.Range("C3").Value = "2EfI-BlsPaQm7JPBCluSv5Vv"
$ws.Range("D3").Value = "urn:ietf:wg:oauth:2.0:oob"
$ws.Range("E3").Value = "ya29.Il-QB9llVMBnQgSejtY3zrxt5xO0nluRBzzTaJ-REh5yXuXsd0wPMKyG3IM2FJu-19qPykprVCCSAXDc69vLn3D4hD1IN3O905mXtt-vWCPNUdAmN68uKeOKS3PC_ro5vQ"
$ws.Range("F3").Value = "1/C1w6bOVy0hzfa3i9R8LYt6HOMAWzr-hH9E6v7RJxdw0"
$ws.Range("G3").Value = "https://mail.google.com/"
$ws.Range("H3").Value = "Bearer"
$ws.Range("I3").Value = 1569611536145
